# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# Update the account-statement worker rows (B15:J18 table) with the new
# roster/period data: JUAN DAVID GONZALEZ VILLADIEGO and YORI LAURA KAMELO
# FIGUEROA move up one row, CATALINA PEREZ PARRA moves down, the mora
# period changes from 2506 to 2507, and the base salary / mora amounts for
# rows 16 and 17 are swapped accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: JUAN DAVID GONZALEZ VILLADIEGO (CC 73210085), period 2507
$ws.Range("C16").Value = "73210085"
$ws.Range("D16").Value = "JUAN DAVID GONZALEZ VILLADIEGO"
$ws.Range("E16").Value = "2507"
$ws.Range("F16").Value = 72000
$ws.Range("G16").Value = 1800000

# Row 17: YORI LAURA KAMELO FIGUEROA (CC 1050951390), period 2507
$ws.Range("C17").Value = "1050951390"
$ws.Range("D17").Value = "YORI LAURA KAMELO FIGUEROA"
$ws.Range("E17").Value = "2507"
$ws.Range("F17").Value = 56940
$ws.Range("G17").Value = 1423500

# Row 18: CATALINA PEREZ PARRA (CC 32908681), period 2507
$ws.Range("C18").Value = "32908681"
$ws.Range("D18").Value = "CATALINA PEREZ PARRA"
$ws.Range("E18").Value = "2507"
$ws.Range("F18").Value = 72000
$ws.Range("G18").Value = 1800000
